$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 30   Number  17"
$ws.Range("C9").Value = "Report Covering the Week  4/24/2023  Through  4/30/2023"

# --- Data table updates (rows 15-30) ---
# Row 15
$ws.Range("N15").Value = -90

# Row 16
$ws.Range("C16").Value = 10
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 150
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = 25
$ws.Range("I16").Value = 80
$ws.Range("J16").Value = 74
$ws.Range("K16").Value = 8.108108108108
$ws.Range("L16").Value = 63.265306122449
$ws.Range("M16").Value = -20.792079207920
$ws.Range("N16").Value = -83.838383838383

# Row 17
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 28
$ws.Range("G17").Value = 33
$ws.Range("H17").Value = -15.151515151515
$ws.Range("I17").Value = 95
$ws.Range("J17").Value = 93
$ws.Range("K17").Value = 2.150537634408
$ws.Range("L17").Value = 11.764705882352
$ws.Range("M17").Value = 58.333333333333
$ws.Range("N17").Value = -65.703971119133

# Row 18
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 133.333333333333
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = 16.666666666666
$ws.Range("I18").Value = 60
$ws.Range("J18").Value = 58
$ws.Range("K18").Value = 3.448275862068
$ws.Range("L18").Value = 50
$ws.Range("M18").Value = -6.25
$ws.Range("N18").Value = -90.839694656488

# Row 19
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -16.666666666666
$ws.Range("F19").Value = 49
$ws.Range("G19").Value = 46
$ws.Range("H19").Value = 6.521739130434
$ws.Range("I19").Value = 170
$ws.Range("J19").Value = 199
$ws.Range("K19").Value = -14.572864321608
$ws.Range("L19").Value = -5.555555555555
$ws.Range("M19").Value = 57.407407407407
$ws.Range("N19").Value = -54.423592493297

# Row 20
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = 16.666666666666
$ws.Range("F20").Value = 23
$ws.Range("G20").Value = 21
$ws.Range("H20").Value = 9.523809523809
$ws.Range("I20").Value = 87
$ws.Range("J20").Value = 125
$ws.Range("K20").Value = -30.4
$ws.Range("L20").Value = 52.631578947368
$ws.Range("M20").Value = 234.615384615385
$ws.Range("N20").Value = -87.606837606837

# Row 21
$ws.Range("C21").Value = 38
$ws.Range("D21").Value = 29
$ws.Range("E21").Value = 31.034482758620
$ws.Range("F21").Value = 129
$ws.Range("G21").Value = 125
$ws.Range("H21").Value = 3.2
$ws.Range("I21").Value = 495
$ws.Range("J21").Value = 558
$ws.Range("K21").Value = -11.290322580645
$ws.Range("L21").Value = 16.470588235294
$ws.Range("M21").Value = 34.510869565217
$ws.Range("N21").Value = -80.648944487881

# Row 22
$ws.Range("C22").Value = 2
$ws.Range("F22").Value = 5
$ws.Range("H22").Value = 66.666666666666
$ws.Range("I22").Value = 16
$ws.Range("K22").Value = 23.076923076923
$ws.Range("L22").Value = 77.777777777777
$ws.Range("M22").Value = 77.777777777777

# Row 23
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "0"
$ws.Range("D22").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0"
$ws.Range("D22").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "***.*"
$ws.Range("E22").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 0
$ws.Range("M23").Value = 10

# Row 24
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = -4.761904761904
$ws.Range("F24").Value = 85
$ws.Range("H24").Value = 7.594936708860
$ws.Range("I24").Value = 403
$ws.Range("J24").Value = 439
$ws.Range("K24").Value = -8.200455580865
$ws.Range("L24").Value = 71.489361702127
$ws.Range("M24").Value = 106.666666666667

# Row 25
$ws.Range("C25").Value = 7
$ws.Range("E25").Value = -12.5
$ws.Range("F25").Value = 38
$ws.Range("G25").Value = 35
$ws.Range("H25").Value = 8.571428571428
$ws.Range("I25").Value = 155
$ws.Range("J25").Value = 147
$ws.Range("K25").Value = 5.442176870748
$ws.Range("L25").Value = 32.478632478632
$ws.Range("M25").Value = -12.921348314606

# Row 26
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = -50
$ws.Range("L26").Value = -53.333333333333

# Row 27
$ws.Range("C27").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D27").Value = 2
$ws.Range("H27").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 66.666666666666
$ws.Range("I27").Value = 20
$ws.Range("J27").Value = 19
$ws.Range("K27").Value = 5.263157894736
$ws.Range("L27").Value = 5.263157894736

# Row 28
$ws.Range("G28").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D28").Value = 1
$ws.Range("H28").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E28").Value = -100
$ws.Range("G28").Value = 2
$ws.Range("J28").Value = 8
$ws.Range("K28").Value = -87.5
$ws.Range("L28").Value = -88.888888888888
$ws.Range("N28").Value = -98.611111111111

# Row 29
$ws.Range("G29").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("D29").Value = 1
$ws.Range("H29").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("E29").Value = -100
$ws.Range("G29").Value = 2
$ws.Range("J29").Value = 8
$ws.Range("K29").Value = -87.5
$ws.Range("L29").Value = -87.5
$ws.Range("N29").Value = -98.507462686567

# Row 30
$ws.Range("K30").Copy()
$ws.Range("L30").PasteSpecial(-4122)
$ws.Range("L30").Value = 100

$excel.CutCopyMode = 0
